$p = $ppt.ActivePresentation

# 1. Bump the version number on the title slide.
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Versie 4.0.0-dev, 06-06-2023"

# 2. Remove slide 20 ("M06: Het project meet kwaliteitsnormen geautomatiseerd en
#    frequent") entirely -- its content overlapped too much with M02, which is
#    now reworded below. Removing this slide shifts all later slides up by one
#    (slide 21 -> 20, slide 22 -> 21, ... slide 33 is gone).
$p.Slides.Item(20).Delete()

# 3. Reword M02's title and body (now slide 8, unaffected by the slide-20
#    deletion since it precedes it).
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "M02: Het project bewaakt continu dat het product aan de kwaliteitsnormen voldoet"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "Projecten bewaken zo snel mogelijk vanaf de start de door het project en ICTU vastgestelde kwaliteitsnormen en voldoen daar zo snel en goed mogelijk aan. De kwaliteit van producten, die nog niet zijn afgerond of nog niet aan de normen voldoen, wordt door het project bewaakt. Het voldoen aan de kwaliteitsnormen is onderdeel van de Definition of Done en herstel van de kwaliteit wordt planmatig opgepakt."
